$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("booking")

# Duplicate the formatting of the row above (row 62) onto the new row 63,
# then overwrite the values/content cell by cell.
$ws.Range("A62:K62").Copy()
$ws.Range("A63:K63").PasteSpecial(-4122)

$ws.Cells.Item(63, 1).Value = "SWAGGER_VALIDATION"
$ws.Cells.Item(63, 2).Value = 5687
$ws.Cells.Item(63, 3).Value = "Samyuktha"
$ws.Cells.Item(63, 4).Value = "Saravanan"
$ws.Cells.Item(63, 5).Value = "'true"
$ws.Cells.Item(63, 6).Value = 46019
$ws.Cells.Item(63, 7).Value = 46022
$ws.Cells.Item(63, 8).Value = "samsaravanan@gmail.com"
$ws.Cells.Item(63, 9).Value = "919710288178"
$ws.Cells.Item(63, 11).Value = "Booking should be updated"
